# Update the "time_taken" timestamps on the "data" sheet (column F, rows 2-15)
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$times = @(
    "2021-10-05 14:18:58.627073",
    "2021-10-05 14:18:58.627080",
    "2021-10-05 14:18:58.627084",
    "2021-10-05 14:18:58.627087",
    "2021-10-05 14:18:58.627090",
    "2021-10-05 14:18:58.627092",
    "2021-10-05 14:18:58.627095",
    "2021-10-05 14:18:58.627097",
    "2021-10-05 14:18:58.627100",
    "2021-10-05 14:18:58.627103",
    "2021-10-05 14:18:58.627105",
    "2021-10-05 14:18:58.627108",
    "2021-10-05 14:18:58.627110",
    "2021-10-05 14:18:58.627113"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $times[$i]
}

# Add a new "metadata" worksheet positioned right after "data"
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

$headerRange = $meta.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Additional findings health related"
$meta.Cells.Item(2, 3).Value = 399
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "0.110"
$meta.Cells.Item(2, 5).Value = "2020-04-21T15:37:36.950369Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:18:58.623044"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/399/?format=json"

$a2 = $meta.Cells.Item(2, 1)
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

# Restore "data" as the active sheet (matches original active tab)
$data.Activate()
